$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.403.70"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.34%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.597.70"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.22%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "588.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.48%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "149.06"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.39%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("E8").Value = "  -1.31%  "

$ws.Range("E9").Value = "  -0.93%  "

$ws.Range("E10").Value = "  -0.63%  "

$ws.Range("E11").Value = "  -0.21%  "

$ws.Range("E12").Value = "  +0.13%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "27.54"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.50%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.064.39"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.30%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "63.237.61"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.39%  "

$ws.Range("E16").Value = "  +3.81%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.597.83"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.14%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.06"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.90%  "

$ws.Range("E19").Value = "  +0.33%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "344.06"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.67%  "

$ws.Range("E21").Value = "  -2.69%  "

$ws.Range("E22").Value = "  +0.16%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.59"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.20%  "

$ws.Range("E24").Value = "  -2.46%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.16"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.39%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.65"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.25%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.26"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.84%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "552.76"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.97%  "

$ws.Range("E29").Value = "  -0.04%  "

$ws.Range("E30").Value = "  -3.40%  "

$ws.Range("E31").Value = "  -2.23%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0₃0851"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.00%  "

$ws.Range("E33").Value = "  -0.35%  "

$ws.Range("E34").Value = "  -0.07%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.04"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.65%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "165.32"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.17%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.411"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.12%  "

$ws.Range("E38").Value = "  -0.03%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.37"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.81%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.92"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.77%  "

$ws.Range("E41").Value = "  -0.06%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "165.13"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.53%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.56%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "22.89"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.73%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0577"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.73%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.10"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.11%  "

$ws.Range("E47").Value = "  +0.28%  "

$ws.Range("E48").Value = "  +0.45%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0958"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.98%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.98"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.66%  "

$ws.Range("E51").Value = "  +11.74%  "
